$d = $word.ActiveDocument
$d.Content.Find.Execute('${no_np}/', $true, $false, $false, $false, $false, $true, 1, $false, '${no_lp}/', 2)
